$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# The two in-flight "running" placeholder timestamp names are now tied to the
# reused SCTv2 corrected pipeline run name (those runs finished/merged in).
$ws.Range("B28").Value = "Pipe_SCTv2_corrected_13-06"
$ws.Range("B30").Value = "Pipe_SCTv2_corrected_13-06"

# New log entries for the Kriegstein SingleR visualization runs.
$ws.Range("A36").Value = "Kriegstein"
$ws.Range("B36").Value = "Pipe_SCTv2_corrected_13-06"
$ws.Range("C36").Value = "SingleR visualization"
$ws.Range("D36").Value = "SCTv2 preSelection aggrFalse labels score.max"

$ws.Range("A37").Value = "Kriegstein"
$ws.Range("B37").Value = "Pipe_SCTv2_corrected_13-06"
$ws.Range("C37").Value = "SingleR visualization"
$ws.Range("D37").Value = "SCTv2 preSelection aggrTrue labels score.max"

$ws.Range("A38").Value = "Kriegstein"
$ws.Range("C38").Value = "SingleR visualization"
$ws.Range("D38").Value = "SCTv2 postSelection aggrFalse labels score.max"

$ws.Range("A39").Value = "Kriegstein"
$ws.Range("C39").Value = "SingleR visualization"
$ws.Range("D39").Value = "SCTv2 postSelection aggrTrue labels score.max"

# These two runs are still active, so their names (timestamps) are filled in
# last, once they finished.
$ws.Range("B38").Value = "2022-06-22 17-28-57"
$ws.Range("B39").Value = "2022-06-22 17-29-21"

# Match the author's final selection/view state.
[void]$ws.Range("B39").Select()
